$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 126
$ws.Range("I2").Value = 318
$ws.Range("J2").Value = 1403
$ws.Range("K2").Value = 8
$ws.Range("L2").Value = 419
$ws.Range("M2").Value = 25
$ws.Range("N2").Value = 265
$ws.Range("P2").Value = 5
$ws.Range("R2").Value = 18
$ws.Range("S2").Value = 134
$ws.Range("T2").Value = 255
$ws.Range("U2").Value = 23
$ws.Range("V2").Value = 2218
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 2185
$ws.Range("Y2").Value = 4
$ws.Range("Z2").Value = 30
$ws.Range("AA2").Value = 9
